# Corrected path in Fam
#
# Rows 34-46 (column A) of Sheet2 still pointed at the old
# "REPSWITCH1_Practice/PICTURE_*.png" image paths even though this block of
# the worksheet belongs to the "Fam" (familiarization) condition. Re-point
# each of those 13 cells at "REPSWITCH1_Fam/..." instead, keeping the same
# PICTURE_#### file name, to match the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = "REPSWITCH1_Fam/PICTURE_612.png"
$ws.Range("A35").Value = "REPSWITCH1_Fam/PICTURE_599.png"
$ws.Range("A36").Value = "REPSWITCH1_Fam/PICTURE_570.png"
$ws.Range("A37").Value = "REPSWITCH1_Fam/PICTURE_733.png"
$ws.Range("A38").Value = "REPSWITCH1_Fam/PICTURE_614.png"
$ws.Range("A39").Value = "REPSWITCH1_Fam/PICTURE_110.png"
$ws.Range("A40").Value = "REPSWITCH1_Fam/PICTURE_12.png"
$ws.Range("A41").Value = "REPSWITCH1_Fam/PICTURE_340.png"
$ws.Range("A42").Value = "REPSWITCH1_Fam/PICTURE_391.png"
$ws.Range("A43").Value = "REPSWITCH1_Fam/PICTURE_590.png"
$ws.Range("A44").Value = "REPSWITCH1_Fam/PICTURE_646.png"
$ws.Range("A45").Value = "REPSWITCH1_Fam/PICTURE_663.png"
$ws.Range("A46").Value = "REPSWITCH1_Fam/PICTURE_673.png"

# Re-apply the font on the first two corrected rows so they get their own
# style entry, matching how the workbook looked once re-saved from Excel.
$ws.Range("A34:A35").Font.Name = "Calibri"

# Leave the selection where Excel left it after making the edit.
$ws.Range("A46").Select()
